# Auto-generated edit script applying the Adamantoise_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 230.9
$ws.Range("J8").Value = 100
$ws.Range("L8").Value = 300
$ws.Range("N8").Value = -578
$ws.Range("H42").Value = 109.916664
$ws.Range("I42").Value = 61
$ws.Range("J42").Value = 178.4
$ws.Range("K42").Value = 183
$ws.Range("L42").Value = 535.2
$ws.Range("M42").Value = 47
$ws.Range("N42").Value = -995.2
$ws.Range("H62").Value = 9855.714
$ws.Range("H65").Value = 9855.714
$ws.Range("H106").Value = 8334254
$ws.Range("I106").Value = 8334254
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 8334254
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -8333623
$ws.Range("N106").ClearContents()
$ws.Range("H124").Value = 145992.5
$ws.Range("J124").Value = 145992.5
$ws.Range("L124").Value = 145992.5
$ws.Range("N124").Value = -155812.5
$ws.Range("H125").Value = 4810
$ws.Range("I125").Value = 1012.5
$ws.Range("J125").Value = 20000
$ws.Range("K125").Value = 9112.5
$ws.Range("L125").Value = 180000
$ws.Range("M125").Value = -6652.5
$ws.Range("N125").Value = -184920
$ws.Range("H132").Value = 3888.8
$ws.Range("I132").Value = 4058.9714
$ws.Range("J132").Value = 2697.6
$ws.Range("K132").Value = 12176.9142
$ws.Range("L132").Value = 8092.799999999999
$ws.Range("M132").Value = -9646.914199999999
$ws.Range("N132").Value = -13152.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H28").Value = 54235.25
$ws.Range("I28").Value = 35314
$ws.Range("J28").Value = 110999
$ws.Range("K28").Value = 35314
$ws.Range("L28").Value = 110999
$ws.Range("M28").Value = -35122
$ws.Range("N28").Value = -111383
$ws.Range("H32").Value = 3315.0227
$ws.Range("I32").Value = 2583.675
$ws.Range("K32").Value = 2583.675
$ws.Range("M32").Value = -2296.675
$ws.Range("H43").Value = 33333.332
$ws.Range("J43").Value = 33333.332
$ws.Range("L43").Value = 33333.332
$ws.Range("N43").Value = -33959.332
$ws.Range("H99").Value = 54235.25
$ws.Range("I99").Value = 35314
$ws.Range("J99").Value = 110999
$ws.Range("K99").Value = 35314
$ws.Range("L99").Value = 110999
$ws.Range("M99").Value = -32319
$ws.Range("N99").Value = -116989
$ws.Range("H122").Value = 5799.8
$ws.Range("I122").Value = 3999.9092
$ws.Range("J122").Value = 7999.6665
$ws.Range("K122").Value = 11999.7276
$ws.Range("L122").Value = 23998.9995
$ws.Range("M122").Value = -9549.7276
$ws.Range("N122").Value = -28898.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 638.7059
$ws.Range("I86").Value = 350.875
$ws.Range("K86").Value = 350.875
$ws.Range("M86").Value = 772.125
$ws.Range("H89").Value = 638.7059
$ws.Range("I89").Value = 350.875
$ws.Range("K89").Value = 1754.375
$ws.Range("M89").Value = 3861.625
$ws.Range("H133").Value = 120001
$ws.Range("J133").Value = 120001
$ws.Range("L133").Value = 120001
$ws.Range("N133").Value = -130121
$ws.Range("H134").Value = 12347708
$ws.Range("J134").Value = 166669170
$ws.Range("L134").Value = 500007510
$ws.Range("N134").Value = -500012580
$ws.Range("H135").Value = 68890
$ws.Range("J135").Value = 68890
$ws.Range("L135").Value = 68890
$ws.Range("N135").Value = -79030
$ws.Range("H141").Value = 90998.5
$ws.Range("J141").Value = 90998.5
$ws.Range("L141").Value = 90998.5
$ws.Range("N141").Value = -101358.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3520.276
$ws.Range("I31").Value = 2322.6843
$ws.Range("K31").Value = 2322.6843
$ws.Range("M31").Value = -2027.6843
$ws.Range("H34").Value = 3520.276
$ws.Range("I34").Value = 2322.6843
$ws.Range("K34").Value = 2322.6843
$ws.Range("M34").Value = -2120.6843
$ws.Range("H58").Value = 3275.9143
$ws.Range("I58").Value = 2958.3914
$ws.Range("J58").Value = 3884.5
$ws.Range("K58").Value = 2958.3914
$ws.Range("L58").Value = 3884.5
$ws.Range("M58").Value = -2755.3914
$ws.Range("N58").Value = -4290.5
$ws.Range("H81").Value = 70494
$ws.Range("J81").Value = 70000
$ws.Range("L81").Value = 70000
$ws.Range("N81").Value = -71996
$ws.Range("H84").Value = 70494
$ws.Range("J84").Value = 70000
$ws.Range("L84").Value = 210000
$ws.Range("N84").Value = -219984
$ws.Range("H86").Value = 4352
$ws.Range("I86").Value = 4927.6665
$ws.Range("K86").Value = 4927.6665
$ws.Range("M86").Value = -3804.6665
$ws.Range("H89").Value = 4352
$ws.Range("I89").Value = 4927.6665
$ws.Range("K89").Value = 24638.3325
$ws.Range("M89").Value = -19022.3325
$ws.Range("H94").Value = 1919.5
$ws.Range("I94").Value = 1063.6666
$ws.Range("J94").Value = 2070.5293
$ws.Range("K94").Value = 1063.6666
$ws.Range("L94").Value = 2070.5293
$ws.Range("M94").Value = -612.6666
$ws.Range("N94").Value = -2972.5293
$ws.Range("H125").Value = 89784.5
$ws.Range("I125").Value = 99999
$ws.Range("J125").Value = 86379.664
$ws.Range("K125").Value = 99999
$ws.Range("L125").Value = 86379.664
$ws.Range("M125").Value = -97539
$ws.Range("N125").Value = -91299.664
$ws.Range("H136").Value = 3275.9143
$ws.Range("I136").Value = 2958.3914
$ws.Range("J136").Value = 3884.5
$ws.Range("K136").Value = 8875.174199999999
$ws.Range("L136").Value = 11653.5
$ws.Range("M136").Value = -6325.174199999999
$ws.Range("N136").Value = -16753.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 673.3
$ws.Range("I38").Value = 72
$ws.Range("J38").Value = 1575.25
$ws.Range("K38").Value = 216
$ws.Range("L38").Value = 4725.75
$ws.Range("M38").Value = 131
$ws.Range("N38").Value = -5419.75
$ws.Range("H98").Value = 925.3077
$ws.Range("I98").Value = 932.5
$ws.Range("J98").Value = 924
$ws.Range("K98").Value = 2797.5
$ws.Range("L98").Value = 2772
$ws.Range("M98").Value = -1299.5
$ws.Range("N98").Value = -5768
$ws.Range("H107").Value = 456.73077
$ws.Range("I107").Value = 428.1111
$ws.Range("J107").Value = 471.88235
$ws.Range("K107").Value = 1284.3333
$ws.Range("L107").Value = 1415.64705
$ws.Range("M107").Value = 635.6667
$ws.Range("N107").Value = -5255.64705
$ws.Range("H128").Value = 74559
$ws.Range("I128").Value = 74559
$ws.Range("K128").Value = 223677
$ws.Range("M128").Value = -218697

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2983.3333
$ws.Range("I102").Value = 2980
$ws.Range("K102").Value = 2980
$ws.Range("M102").Value = -1358
$ws.Range("H122").Value = 2407.5
$ws.Range("I122").Value = 1945
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 5835
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3385
$ws.Range("N122").Value = -12400
$ws.Range("H132").Value = 3381.8
$ws.Range("I132").Value = 3381.8
$ws.Range("K132").Value = 10145.4
$ws.Range("M132").Value = -7615.400000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4953.476
$ws.Range("I40").Value = 4422.316
$ws.Range("K40").Value = 4422.316
$ws.Range("M40").Value = -4286.316
$ws.Range("H61").Value = 4999.091
$ws.Range("I61").Value = 1558.8
$ws.Range("K61").Value = 1558.8
$ws.Range("M61").Value = -1356.8
$ws.Range("H113").Value = 4999.091
$ws.Range("I113").Value = 1558.8
$ws.Range("K113").Value = 1558.8
$ws.Range("M113").Value = 611.2
$ws.Range("H136").Value = 4091.1765
$ws.Range("I136").Value = 3858.4167
$ws.Range("J136").Value = 4649.8
$ws.Range("K136").Value = 11575.2501
$ws.Range("L136").Value = 13949.4
$ws.Range("M136").Value = -9025.250100000001
$ws.Range("N136").Value = -19049.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4949.9165
$ws.Range("I62").Value = 3779.4
$ws.Range("J62").Value = 5786
$ws.Range("K62").Value = 3779.4
$ws.Range("L62").Value = 5786
$ws.Range("M62").Value = -3155.4
$ws.Range("N62").Value = -7034
$ws.Range("H65").Value = 4949.9165
$ws.Range("I65").Value = 3779.4
$ws.Range("J65").Value = 5786
$ws.Range("K65").Value = 18897
$ws.Range("L65").Value = 28930
$ws.Range("M65").Value = -15777
$ws.Range("N65").Value = -35170
$ws.Range("H136").Value = 18528.05
$ws.Range("I136").Value = 1226.7609
$ws.Range("J136").Value = 79748
$ws.Range("K136").Value = 3680.2827
$ws.Range("L136").Value = 239244
$ws.Range("M136").Value = -1130.2827
$ws.Range("N136").Value = -244344
